$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = 5
$ws.Range("E17").Value = 14
$ws.Range("E19").Value = 1
